$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two trailing unused template rows (41 and 42) first, from
# bottom to top so row numbers above them are not disturbed prematurely.
$ws.Rows.Item(42).Delete() | Out-Null
$ws.Rows.Item(41).Delete() | Out-Null

# Fill in row 39 with the new LeetCode entry (previously a blank "Anna" template row).
$ws.Range("B39").Value = "Stephan"
$ws.Range("D39").Value = "453. Minimum Moves to Equal Array Elements"
$ws.Range("E39").Value = "2021/01/21"
$ws.Range("F39").Value = "Array, Large Integer, Sort"
$ws.Range("F39").WrapText = $true
$ws.Range("G39").Value = "Completed"
$ws.Rows.Item(39).RowHeight = 28

# Update the visible window / selection to match the new scroll position.
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Range("D39:G39").Select() | Out-Null
